$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata_input")

# Rename the "Time" column header to "Timepoint"
$ws.Range("B1").Value = "Timepoint"
